$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2184343333333333
$ws.Range("H2").Value = 0.655303
$ws.Range("I2").Value = 0.008416673064019609
$ws.Range("J2").Value = 0.00841667306401961
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03887266666666667
$ws.Range("N2").Value = 0.116618
$ws.Range("Q2").Value = 0.008491125028222222
$ws.Range("R2").Value = 0.07642012525399999
$ws.Range("S2").Value = 0.008416673064019609
$ws.Range("T2").Value = 0.00841667306401961

# Row 3
$ws.Range("I3").Value = 0.1618270290283213
$ws.Range("J3").Value = 0.1618270290283213
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03887266666666667
$ws.Range("N3").Value = 0.116618
$ws.Range("Q3").Value = 0.1632585139013333
$ws.Range("R3").Value = 1.469326625112
$ws.Range("S3").Value = 0.1618270290283213
$ws.Range("T3").Value = 0.1618270290283213

# Row 4
$ws.Range("G4").Value = 3.307112333333333
$ws.Range("H4").Value = 9.921336999999999
$ws.Range("I4").Value = 0.1274290669918512
$ws.Range("J4").Value = 0.1274290669918513
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03887266666666667
$ws.Range("N4").Value = 0.116618
$ws.Range("Q4").Value = 0.1285562753628889
$ws.Range("R4").Value = 1.157006478266
$ws.Range("S4").Value = 0.1274290669918512
$ws.Range("T4").Value = 0.1274290669918513

# Row 5
$ws.Range("G5").Value = 18.22719966666667
$ws.Range("H5").Value = 54.681599
$ws.Range("I5").Value = 0.7023272309158078
$ws.Range("J5").Value = 0.7023272309158078
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03887266666666667
$ws.Range("N5").Value = 0.116618
$ws.Range("Q5").Value = 0.7085398569091111
$ws.Range("R5").Value = 6.376858712182
$ws.Range("S5").Value = 0.7023272309158078
$ws.Range("T5").Value = 0.7023272309158078
